# Update daily and intraday volume tables for Bond and E-mini Futures
# - Adjusted average daily volume and differences for 10yr Bond Futures in Table 1.
# - Modified average daily volume and differences for E-mini Futures in Table 1 and Table 2.
# - Updated observations count for E-mini Futures in Table 1 and Table 2.
# - Revised average daily volume and differences for E-mini Futures in Table 3 (Post-ZLB, Pre-ZLB, ZLB).
# - Corrected observations count for E-mini Futures in Table 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small last-digit roundoff refreshes scattered through the table (rows 2,3,5,8,9,11,14,15) ---
$ws.Range("E2").Value   = 145.5013476064509
$ws.Range("V2").Value   = 60.89029304029305
$ws.Range("D3").Value   = 50.68269230769231
$ws.Range("J3").Value   = 38.14438839848676
$ws.Range("V3").Value   = 50.88207417582417
$ws.Range("K5").Value   = 155.1003991844884
$ws.Range("Q5").Value   = 249.4473527026094
$ws.Range("J8").Value   = 143.4974779319041
$ws.Range("AB8").Value  = 41.7595134032634
$ws.Range("D9").Value   = 94.74999999999999
$ws.Range("D11").Value  = 181.5477667493796
$ws.Range("P11").Value  = 224.7522250476797
$ws.Range("K14").Value  = 243.1573772934979
$ws.Range("P14").Value  = 231.9639224411952
$ws.Range("AB15").Value = 18.90740093240094

# --- row 26: Panel B (E-mini Futures) - "Ann Window Volume" stats, recomputed with revised obs count ---
$ws.Range("D26").Value  = 4001.254480286738
$ws.Range("E26").Value  = 1595.588078664807
$ws.Range("F26").Value  = 2630.258064516129
$ws.Range("G26").Value  = 3908.290322580645
$ws.Range("H26").Value  = 5093.16129032258
$ws.Range("I26").Value  = 63
$ws.Range("J26").Value  = 3338.169138693728
$ws.Range("K26").Value  = 1224.403855081451
$ws.Range("L26").Value  = 2310.696721311475
$ws.Range("M26").Value  = 3172.344262295082
$ws.Range("N26").Value  = 4256.081967213115
$ws.Range("O26").Value  = 63
$ws.Range("P26").Value  = 3423.874196510561
$ws.Range("Q26").Value  = 1213.705166738312
$ws.Range("R26").Value  = 2438.256198347108
$ws.Range("S26").Value  = 3503.066115702479
$ws.Range("T26").Value  = 4219.285123966942
$ws.Range("U26").Value  = 63
$ws.Range("V26").Value  = 3057.257218442932
$ws.Range("W26").Value  = 1012.369413239874
$ws.Range("X26").Value  = 2292.154761904762
$ws.Range("Y26").Value  = 2945.295238095238
$ws.Range("Z26").Value  = 3761.464285714286
$ws.Range("AA26").Value = 63
$ws.Range("AB26").Value = 876.6249278499279
$ws.Range("AC26").Value = 232.2364393614542
$ws.Range("AG26").Value = 63

# --- row 27: Panel B (E-mini Futures) - "Diff (Ann - Non)" ---
$ws.Range("D27").Value  = 1903.087685611879
$ws.Range("J27").Value  = 1236.004553734062
$ws.Range("P27").Value  = 1285.900400104946
$ws.Range("V27").Value  = 800.25589569161
$ws.Range("AB27").Value = 28.82795214045213

# --- row 28: Panel B (E-mini Futures) - "# Obs" ---
$ws.Range("D28").Value  = 63
$ws.Range("J28").Value  = 63
$ws.Range("P28").Value  = 63
$ws.Range("V28").Value  = 63
$ws.Range("AB28").Value = 63

Write-Output "Applied 52 cell updates to Sheet1 (Table3 Post-ZLB intraday volume)."
